$d = $word.ActiveDocument

# --- First paragraph: replace the topic placeholder text (and drop the
#     trailing space run that followed it) ---
$d.Content.Find.Execute("**ID__AFFARS_5336_topic_15__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_5336_609_1__ID**", 2)

# --- First paragraph: paragraph border (space-only, no visible lines) and
#     updated left indent ---
$p1 = $d.Paragraphs(1)
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25
